$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Plan")

# Update the "Player Help" acceptance-criteria rows (47-49): the old
# "redirected to a help page" / "return to my game" criteria plus their
# "KC; 4/11/18" comments and "Pass" statuses are replaced with the new,
# reset criteria text and cleared status/comment cells.
$ws.Range("B48").Value = "Given I am a player when I click the request help link then I expect to be shown where I can make a move."
$ws.Range("B49").Value = "Given I am a player when I make a move after requesting help then I expect the valid move spaces to no longer be indicated."

# Clear out the Sprint status (column G) for all three rows.
$ws.Range("G47").ClearContents()
$ws.Range("G48").ClearContents()
$ws.Range("G49").ClearContents()

# Clear out the tester comment cells (column H) entirely.
$ws.Range("H47").Clear()
$ws.Range("H48").Clear()
$ws.Range("H49").Clear()

# Shrink the row heights now that the rows hold shorter text.
$ws.Rows.Item(47).RowHeight = 25.35
$ws.Rows.Item(48).RowHeight = 25.35
$ws.Rows.Item(49).RowHeight = 25.35

# Update the window scroll/selection to match the reset view.
[void]$ws.Activate()
[void]$ws.Range("A34").Select()
[void]$ws.Range("G49").Select()
